# Updated FIN model - 2025-08-25 16:59
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet (was "Sheet1") to "misc."
$ws.Name = "misc."

# --- New header cells on row 4 (K4/L4), styled like the rest of the header row ---
$ws.Range("K4").Value = "other_indexes"
$ws.Range("L4").Value = "commodity"
$ws.Range("K4:L4").Style = "Heading 3"

# --- New data row 11 (flo_emis / gas / ccs exemption for co2captured->co2) ---
$ws.Range("B11").Value = "flo_emis"
$ws.Range("D11").Value = "gas"
$ws.Range("L11").Value = "co2captured"
$ws.Range("K11").Value = "co2"
$ws.Range("E11").Value = "*ccs,*ccs-rf"
$ws.Range("H11").Value = 0.95

# --- New data row 12 (flo_emis / coal,oil / ccs exemption for co2captured->co2) ---
$ws.Range("B12").Value = "flo_emis"
$ws.Range("D12").Value = "coal,oil"
$ws.Range("E12").Value = "*ccs,*ccs-rf"
$ws.Range("H12").Value = 0.85
$ws.Range("K12").Value = "co2"
$ws.Range("L12").Value = "co2captured"

# --- Column width touch-ups (column E widened for the new long text, column K newly sized) ---
$ws.Columns.Item(5).ColumnWidth = 9.498697916666666
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666

# --- Selection left where the author's cursor ended up ---
$ws.Range("D13").Select()
